# Feature: add arrows (arrow_n).
#
# The "meta" sheet (first worksheet) holds key/value metadata rows in
# columns A (key, bold/orange style) and B (value). This adds a new
# "style" / "default" metadata row right where the old trailing blank
# templated row (A6, style-only) used to live, and pushes that blank
# templated row down to row 7 - matching how the sheet is re-generated
# with one extra metadata entry.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert a new row at 6, pushing the old (blank, style-only) row 6 down to
# row 7. The inserted row inherits the formatting of row 6 (A gets s="1").
$ws1.Rows(6).Insert()

# Populate the newly freed row 6 with the new metadata pair.
$ws1.Range("A6").Value = "style"
$ws1.Range("B6").Value = "default"
